$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove column M entirely; this shifts column N (and beyond) one column to
# the left, so the former column N values become the new column M values.
$ws.Columns.Item(13).Delete()
